# "grace marks form model"
# Bump the Sessional Exam / End Term marks for the first two students so the
# Overall Mark (C+D) reflects grace-marks adjustments. The Overall Mark column
# is a formula (C+D) so it recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Roll No 2023-133): Sessional Exam 20 -> 100, End Term 50 -> 300
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 300

# Row 3 (Roll No 2023-156): Sessional Exam 18 -> 150, End Term 45 -> 340
$ws.Range("C3").Value = 150
$ws.Range("D3").Value = 340

# Leave the active selection where the editor last left off.
$ws.Range("D4").Select() | Out-Null
